$d = $word.ActiveDocument

$p2 = $d.Paragraphs(2).Range
$xml2 = @'
<w:p w:rsidR="00831858" w:rsidRDefault="00831858">
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Ffhq</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
    </w:p>
'@
$p2.InsertXML($xml2)

$p3 = $d.Paragraphs(3).Range
$xml3 = @'
<w:p w:rsidR="00831858" w:rsidRDefault="00831858">
      <w:r>
        <w:t>Celeb</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">-&gt; </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>atttrative</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> attribute</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">이 </w:t>
      </w:r>
    </w:p>
'@
$p3.InsertXML($xml3)

$p4 = $d.Paragraphs(4).Range
$xml4 = @'
<w:p w:rsidR="00831858" w:rsidRDefault="00DD2A26">
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">선하다 </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">/ </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">악하다 </w:t>
      </w:r>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>외향</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> 내향 </w:t>
      </w:r>
    </w:p>
'@
$p4.InsertXML($xml4)

$p8 = $d.Paragraphs(8).Range
$xml8 = @'
<w:p w:rsidR="006E55D9" w:rsidRDefault="006E55D9" w:rsidP="006E55D9">
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:ind w:leftChars="0"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>2개:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">선함 </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">vs </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">악함 </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">/ </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>외향</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">vs </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">내향 </w:t>
      </w:r>
    </w:p>
'@
$p8.InsertXML($xml8)

$p10 = $d.Paragraphs(10).Range
$xml10 = @'
<w:p w:rsidR="00034933" w:rsidRDefault="00034933" w:rsidP="00034933">
      <w:pPr>
        <w:ind w:left="96"/>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Ffhq</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>: 8</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">만장 </w:t>
      </w:r>
      <w:r>
        <w:t>-&gt; filtering</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>을 해서</w:t>
      </w:r>
      <w:r w:rsidR="00AF61DF">
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r w:rsidR="003E61C7">
        <w:t>2000</w:t>
      </w:r>
      <w:r w:rsidR="00AF61DF">
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>개 중</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>600</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">장 정도 나옴 </w:t>
      </w:r>
    </w:p>
'@
$p10.InsertXML($xml10)

$p15 = $d.Paragraphs(15).Range
$xml15 = @'
<w:p w:rsidR="00A97072" w:rsidRDefault="00A97072" w:rsidP="00A97072">
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:ind w:leftChars="0"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>6</w:t>
      </w:r>
      <w:r>
        <w:t>00</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">장에 대해 </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">-&gt; 50% </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">정도 </w:t>
      </w:r>
      <w:r>
        <w:t>consensus : 300</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">장 </w:t>
      </w:r>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t>만장일치거나</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">명이 무효 </w:t>
      </w:r>
      <w:r>
        <w:t>-&gt; strict)</w:t>
      </w:r>
    </w:p>
'@
$p15.InsertXML($xml15)

$rCombined = $d.Paragraphs(16).Range
$xmlCombined = @'
<w:p w:rsidR="00A97072" w:rsidRDefault="00A97072" w:rsidP="00A97072">
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:ind w:leftChars="0"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">한명이 </w:t>
      </w:r>
      <w:r>
        <w:t>inverse</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
        </w:rPr>
        <w:t xml:space="preserve">를 던지면 연한 </w:t>
      </w:r>
      <w:r>
        <w:t>threshold</w:t>
      </w:r>
    </w:p>
<w:p>
      <w:pPr>
        <w:pStyle w:val="a3"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:ind w:leftChars="0"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Keyword는 attractive / extrovert / </w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
'@
$rCombined.InsertXML($xmlCombined)
